# New crime data collected - weekly CompStat update for the 108th Precinct.
# Bumps the report volume/number + covering-week dates, and refreshes the
# Crime Complaints grid (rows 14-27) with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text: volume/number and the covering-week date range ---
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# --- Row 14: MURDER ---
$ws.Range("N14").Value = -90

# --- Row 16: ROBBERY ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 14
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 183
$ws.Range("J16").Value = 119
$ws.Range("K16").Value = 53.781512605042
$ws.Range("L16").Value = 63.392857142857
$ws.Range("M16").Value = -1.081081081081
$ws.Range("N16").Value = -81.114551083591

# --- Row 17: FEL. ASSAULT ---
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -42.857142857142
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -26.315789473684
$ws.Range("I17").Value = 210
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 5
$ws.Range("L17").Value = 39.072847682119
$ws.Range("M17").Value = 72.131147540983
$ws.Range("N17").Value = -30.232558139534

# --- Row 18: BURGLARY ---
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 163
$ws.Range("J18").Value = 174
$ws.Range("K18").Value = -6.321839080459
$ws.Range("L18").Value = -4.678362573099
$ws.Range("M18").Value = -33.469387755102
$ws.Range("N18").Value = -88.964116452268

# --- Row 19: GR. LARCENY ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -36.842105263157
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = -39.130434782608
$ws.Range("I19").Value = 650
$ws.Range("J19").Value = 457
$ws.Range("K19").Value = 42.231947483588
$ws.Range("L19").Value = 63.727959697733
$ws.Range("M19").Value = 47.058823529411
$ws.Range("N19").Value = -25.968109339407

# --- Row 20: G.L.A. ---
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 600
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 36.363636363636
$ws.Range("I20").Value = 223
$ws.Range("J20").Value = 171
$ws.Range("K20").Value = 30.409356725146
$ws.Range("L20").Value = 39.375
$ws.Range("M20").Value = 16.753926701570
$ws.Range("N20").Value = -88.805220883534

# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -18.181818181818
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 126
$ws.Range("H21").Value = -19.841269841269
$ws.Range("I21").Value = 1449
$ws.Range("J21").Value = 1135
$ws.Range("K21").Value = 27.665198237885
$ws.Range("L21").Value = 44.322709163346
$ws.Range("M21").Value = 20.448877805486
$ws.Range("N21").Value = -74.303954601879

# --- Row 22: TRANSIT --- (C22 flips from a text "0" to a real number, so
# restyle it to the same numeric format the rest of the row uses)
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -75
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 76
$ws.Range("J22").Value = 44
$ws.Range("K22").Value = 72.727272727272
$ws.Range("L22").Value = 145.161290322581
$ws.Range("M22").Value = 55.102040816326

# --- Row 24: PETIT LARCENY ---
$ws.Range("C24").Value = 46
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 109.090909090909
$ws.Range("F24").Value = 157
$ws.Range("G24").Value = 110
$ws.Range("H24").Value = 42.727272727272
$ws.Range("I24").Value = 1471
$ws.Range("J24").Value = 1222
$ws.Range("K24").Value = 20.376432078559
$ws.Range("L24").Value = 69.080459770114
$ws.Range("M24").Value = 64.725643896976

# --- Row 25: MISD. ASSAULT ---
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 24.242424242424
$ws.Range("I25").Value = 535
$ws.Range("J25").Value = 449
$ws.Range("K25").Value = 19.153674832962
$ws.Range("L25").Value = 69.841269841269
$ws.Range("M25").Value = 11.924686192468

# --- Row 26: UCR RAPE* ---
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0

# --- Row 27: OTHER SEX CRIMES ---
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 37.5
$ws.Range("I27").Value = 91
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = 44.444444444444
$ws.Range("L27").Value = 49.180327868852
